$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "release/8.0.6"
$ws.Range("B9").Value = "X"
$ws.Range("C9").Value = "X"
$ws.Range("D9").Value = "X"
$ws.Range("E9").Value = "X"
